$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value2 = "70.071.76"
$cell.Style = "Normal"
$ws.Range("E2").Value2 = "  +4.94%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value2 = "3.601.09"
$cell.Style = "Normal"
$ws.Range("E3").Value2 = "  +4.94%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value2 = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value2 = "  +0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value2 = "589.02"
$cell.Style = "Normal"
$ws.Range("E5").Value2 = "  +3.45%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value2 = "190.19"
$cell.Style = "Normal"
$ws.Range("E6").Value2 = "  +3.78%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value2 = "0.645"
$cell.Style = "Normal"
$ws.Range("E7").Value2 = "  +1.76%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value2 = "3.592.18"
$cell.Style = "Normal"
$ws.Range("E8").Value2 = "  +4.94%  "

$ws.Range("E9").Value2 = "  +0.04%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value2 = "0.178"
$cell.Style = "Normal"
$ws.Range("E10").Value2 = "  -0.24%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value2 = "0.661"
$cell.Style = "Normal"
$ws.Range("E11").Value2 = "  +2.56%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value2 = "57.90"
$cell.Style = "Normal"
$ws.Range("E12").Value2 = "  +4.89%  "

$ws.Range("E13").Value2 = "  +3.18%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value2 = "9.75"
$cell.Style = "Normal"
$ws.Range("E14").Value2 = "  +4.11%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value2 = "4.167.41"
$cell.Style = "Normal"
$ws.Range("E15").Value2 = "  +5.05%  "

$ws.Range("B16").Value2 = "Chainlink"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value2 = "19.35"
$cell.Style = "Normal"
$ws.Range("E16").Value2 = "  +4.77%  "

$ws.Range("B17").Value2 = "WrappedEther"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value2 = "3.591.57"
$cell.Style = "Normal"
$ws.Range("E17").Value2 = "  +4.68%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value2 = "69.961.54"
$cell.Style = "Normal"
$ws.Range("E18").Value2 = "  +4.88%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value2 = "12.46"
$cell.Style = "Normal"
$ws.Range("E19").Value2 = "  +3.37%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value2 = "0.121"
$cell.Style = "Normal"
$ws.Range("E20").Value2 = "  +0.24%  "

$ws.Range("E21").Value2 = "  +3.95%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value2 = "494.91"
$cell.Style = "Normal"
$ws.Range("E22").Value2 = "  +5.15%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value2 = "17.43"
$cell.Style = "Normal"
$ws.Range("E23").Value2 = "  +18.33%  "

$ws.Range("E24").Value2 = "  +7.31%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value2 = "4.46"
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value2 = "90.70"
$cell.Style = "Normal"
$ws.Range("E26").Value2 = "  +0.95%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value2 = "3.12"
$cell.Style = "Normal"
$ws.Range("E27").Value2 = "  +5.65%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value2 = "11.08"
$cell.Style = "Normal"
$ws.Range("E28").Value2 = "  +1.45%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value2 = "9.41"
$cell.Style = "Normal"
$ws.Range("E29").Value2 = "  +5.52%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value2 = "32.24"
$cell.Style = "Normal"
$ws.Range("E30").Value2 = "  +2.29%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value2 = "7.63"
$cell.Style = "Normal"
$ws.Range("E31").Value2 = "  +9.80%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value2 = "12.22"
$cell.Style = "Normal"
$ws.Range("E32").Value2 = "  +5.15%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value2 = "620.41"
$cell.Style = "Normal"
$ws.Range("E33").Value2 = "  +5.45%  "

$ws.Range("E34").Value2 = "  +7.06%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value2 = "65.04"
$cell.Style = "Normal"
$ws.Range("E35").Value2 = "  +3.81%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value2 = "0.0₃0818"
$cell.Style = "Normal"
$ws.Range("E36").Value2 = "  +6.93%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value2 = "38.13"
$cell.Style = "Normal"
$ws.Range("E37").Value2 = "  +4.42%  "

$ws.Range("B38").Value2 = "TheGraph"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value2 = "0.404"
$cell.Style = "Normal"
$ws.Range("E38").Value2 = "  +4.60%  "

$ws.Range("B39").Value2 = "Dai"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value2 = "1.00"
$cell.Style = "Normal"
$ws.Range("E39").Value2 = "  +0.03%  "

$ws.Range("E40").Value2 = "  -0.32%  "

$ws.Range("E41").Value2 = "  +0.57%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value2 = "3.305.79"
$cell.Style = "Normal"
$ws.Range("E42").Value2 = "  +5.52%  "

$ws.Range("E43").Value2 = "  +5.10%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value2 = "0.0446"
$cell.Style = "Normal"
$ws.Range("E44").Value2 = "  +4.46%  "

$ws.Range("E45").Value2 = "  +4.98%  "

$ws.Range("B46").Value2 = "Stellar"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value2 = "0.137"
$cell.Style = "Normal"
$ws.Range("E46").Value2 = "  +1.86%  "

$ws.Range("B47").Value2 = "ApeXProtocol"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value2 = "3.25"
$cell.Style = "Normal"
$ws.Range("E47").Value2 = "  +1.89%  "

$ws.Range("E48").Value2 = "  +5.39%  "

$ws.Range("B49").Value2 = "dogwifhat"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value2 = "2.70"
$cell.Style = "Normal"
$ws.Range("E49").Value2 = "  -4.20%  "

$ws.Range("B50").Value2 = "LidoDAOToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value2 = "3.31"
$cell.Style = "Normal"
$ws.Range("E50").Value2 = "  +5.37%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value2 = "0.998"
$cell.Style = "Normal"
$ws.Range("E51").Value2 = "  -0.03%  "
